# India Super League - atualização de bases (29-02-2024 07:50)
# Fills in match results (and derived P/L columns) for rows 97 and 98,
# updates the odds movement for row 99, and appends a new fixture as row 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 97 : Hyderabad FC vs Punjab FC -> result filled in (FTHG/FTAG/FTR) plus
# refreshed closing odds and P/L figures.
# ---------------------------------------------------------------------------
$ws.Cells.Item(97,8).Value  = 0        # H97  FTHG
$ws.Cells.Item(97,9).Value  = 2        # I97  FTAG
$ws.Cells.Item(97,10).Value = "A"      # J97  FTR

$ws.Cells.Item(97,14).Value = 4.75     # N97  oddH
$ws.Cells.Item(97,15).Value = 3.75     # O97  oddD
$ws.Cells.Item(97,16).Value = 1.55     # P97  oddA

$ws.Cells.Item(97,18).Value = 2.05     # R97  oddAHH
$ws.Cells.Item(97,19).Value = 1.8      # S97  oddAHA

$ws.Cells.Item(97,21).Value = 1.85     # U97  oddAHOver
$ws.Cells.Item(97,22).Value = 2        # V97  oddAHUnder

$ws.Cells.Item(97,23).Value = -1       # W97  PLH
$ws.Cells.Item(97,24).Value = -1       # X97  PLD
$ws.Cells.Item(97,25).Value = 0.55     # Y97  PLA
$ws.Cells.Item(97,26).Value = -1       # Z97  PL_Ahh
$ws.Cells.Item(97,27).Value = 0.8      # AA97 PL_Aha
$ws.Cells.Item(97,28).Value = -1       # AB97 PL_AhOver
$ws.Cells.Item(97,29).Value = 1        # AC97 PL_AhUnder

# ---------------------------------------------------------------------------
# Row 98 : Mumbai City FC vs FC Goa -> result filled in plus refreshed odds.
# ---------------------------------------------------------------------------
$ws.Cells.Item(98,8).Value  = 1        # H98  FTHG
$ws.Cells.Item(98,9).Value  = 1        # I98  FTAG
$ws.Cells.Item(98,10).Value = "D"      # J98  FTR

$ws.Cells.Item(98,14).Value = 2.1      # N98  oddH
$ws.Cells.Item(98,16).Value = 3.1      # P98  oddA
$ws.Cells.Item(98,17).Value = -0.25    # Q98  Ah

$ws.Cells.Item(98,18).Value = 1.875    # R98  oddAHH
$ws.Cells.Item(98,19).Value = 1.975    # S98  oddAHA

$ws.Cells.Item(98,21).Value = 1.925    # U98  oddAHOver
$ws.Cells.Item(98,22).Value = 1.925    # V98  oddAHUnder

$ws.Cells.Item(98,23).Value = -1       # W98  PLH
$ws.Cells.Item(98,24).Value = 2.4      # X98  PLD
$ws.Cells.Item(98,25).Value = -1       # Y98  PLA
$ws.Cells.Item(98,26).Value = -0.5     # Z98  PL_Ahh
$ws.Cells.Item(98,27).Value = 0.4875   # AA98 PL_Aha
$ws.Cells.Item(98,28).Value = -1       # AB98 PL_AhOver
$ws.Cells.Item(98,29).Value = 0.925    # AC98 PL_AhUnder

# ---------------------------------------------------------------------------
# Row 99 : Odisha FC vs Bengaluru -> odds movement refreshed (no result yet).
# ---------------------------------------------------------------------------
$ws.Cells.Item(99,18).Value = 1.775    # R99  oddAHH
$ws.Cells.Item(99,19).Value = 2.025    # S99  oddAHA
$ws.Cells.Item(99,21).Value = 1.825    # U99  oddAHOver
$ws.Cells.Item(99,22).Value = 1.975    # V99  oddAHUnder

# ---------------------------------------------------------------------------
# Row 100 : new fixture appended - Mohun Bagan SG vs Jamshedpur FC.
# Copy row 99 so the bold/centered id style (A) and the date style (E) carry
# across, then overwrite every cell with the new match's data.
# ---------------------------------------------------------------------------
$ws.Range("A99").Copy($ws.Range("A100"))
$ws.Range("E99").Copy($ws.Range("E100"))

$ws.Cells.Item(100,1).Value  = 98                  # A100 id
$ws.Cells.Item(100,2).Value  = 7749759             # B100 id (source)
$ws.Cells.Item(100,3).Value  = "India Super League" # C100 Div
$ws.Cells.Item(100,4).Value  = "India Super League" # D100 Div Original Name
$ws.Cells.Item(100,5).Value  = 45352.45833333334   # E100 Date
$ws.Cells.Item(100,6).Value  = "Mohun Bagan SG"    # F100 HomeTeam
$ws.Cells.Item(100,7).Value  = "Jamshedpur FC"     # G100 AwayTeam

$ws.Cells.Item(100,11).Value = 1.615   # K100 oddH_op
$ws.Cells.Item(100,12).Value = 3.75    # L100 oddD_op
$ws.Cells.Item(100,13).Value = 5       # M100 oddA_op
$ws.Cells.Item(100,14).Value = 1.615   # N100 oddH
$ws.Cells.Item(100,15).Value = 3.75    # O100 oddD
$ws.Cells.Item(100,16).Value = 5       # P100 oddA
$ws.Cells.Item(100,17).Value = -0.75   # Q100 Ah
$ws.Cells.Item(100,18).Value = 1.825   # R100 oddAHH
$ws.Cells.Item(100,19).Value = 1.975   # S100 oddAHA
$ws.Cells.Item(100,20).Value = 2.5     # T100 AhOU
$ws.Cells.Item(100,21).Value = 1.8     # U100 oddAHOver
$ws.Cells.Item(100,22).Value = 2       # V100 oddAHUnder
$ws.Cells.Item(100,23).Value = 0       # W100 PLH
$ws.Cells.Item(100,24).Value = 0       # X100 PLD
$ws.Cells.Item(100,25).Value = 0       # Y100 PLA
$ws.Cells.Item(100,26).Value = 0       # Z100 PL_Ahh
$ws.Cells.Item(100,27).Value = 0       # AA100 PL_Aha
